$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing "mobile" numbers (A2:A4) are stored as text (shared strings),
# so the new entry should follow the same convention rather than being
# auto-converted to a numeric value. Temporarily mark the cell as Text so
# the value is kept as a string, then restore the default style so the
# cell's formatting matches its neighbours.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "7383211888"
$ws.Range("A5").Style = "Normal"
